$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($rng, [string]$newName) {
    $shp = $rng.InlineShapes.Item(1)
    # Renaming the picture in a header/footer range directly can hit a
    # stale-handle error in this host; selecting it first and renaming
    # through the Selection's InlineShapes collection works reliably.
    $shp.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Headers: BTec_Logo-Orange, currently "image1.jpg" -> "image2.jpg"
Rename-InlinePicture $sec.Headers.Item(1).Range "image2.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"

# Footers: PearsonLogo, currently "image2.png" -> "image1.png"
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"
